# Refresh crypto market snapshot (rank-order swaps, prices, and 1h volume deltas)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.416.51"
$ws.Range("E2").Value = "  +1.60%  "

# Row 3
$ws.Range("D3").Value = "1.906.68"
$ws.Range("E3").Value = "  +0.11%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'326.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "

# Row 6
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").Value = "'0.4667"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.54%  "

# Row 8
$ws.Range("D8").Value = "'0.4078"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.68%  "

# Row 9
$ws.Range("D9").Value = "'47.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.44%  "

# Row 10
$ws.Range("E10").Value = "  +0.16%  "

# Row 11
$ws.Range("E11").Value = "  +0.28%  "

# Row 12
$ws.Range("D12").Value = "'22.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.03%  "

# Row 13
$ws.Range("D13").Value = "1.929.15"
$ws.Range("E13").Value = "  +1.06%  "

# Row 14
$ws.Range("D14").Value = "'5.932"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.19%  "

# Row 15
$ws.Range("D15").Value = "'7.109"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "

# Row 16
$ws.Range("D16").Value = "'89.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

# Row 17
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$ws.Range("D18").Value = "'0.06595"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "

# Row 19
$ws.Range("D19").Value = "'0.00001025"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "

# Row 20
$ws.Range("D20").Value = "'17.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.64%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").Value = "29.429.13"
$ws.Range("E22").Value = "  +1.59%  "

# Row 23
$ws.Range("D23").Value = "'5.521"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "

# Row 24
$ws.Range("D24").Value = "'11.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.66%  "

# Row 25
$ws.Range("D25").Value = "'2.211"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.19%  "

# Row 26
$ws.Range("D26").Value = "2.121.09"
$ws.Range("E26").Value = "  -0.51%  "

# Row 27
$ws.Range("D27").Value = "'153.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.81%  "

# Row 28
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("D29").Value = "'2.129"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.62%  "

# Row 30
$ws.Range("D30").Value = "'5.698"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.63%  "

# Row 31
$ws.Range("D31").Value = "'116.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.74%  "

# Row 32
$ws.Range("D32").Value = "'1.070"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.26%  "

# Row 33
$ws.Range("D33").Value = "'0.09485"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.07%  "

# Row 34
$ws.Range("D34").Value = "'1.417"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "

# Row 35
$ws.Range("E35").Value = "  -0.56%  "

# Row 36
$ws.Range("D36").Value = "'5.374"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.66%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02253"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.32%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06072"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "

# Row 39
$ws.Range("D39").Value = "'8.337"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "

# Row 40
$ws.Range("D40").Value = "'1.173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "

# Row 41
$ws.Range("D41").Value = "'0.5858"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "

# Row 42
$ws.Range("E42").Value = "  +0.63%  "

# Row 43
$ws.Range("E43").Value = "  -0.33%  "

# Row 44
$ws.Range("D44").Value = "'1.300"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.00%  "

# Row 45
$ws.Range("D45").Value = "'2.426"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.14%  "

# Row 46
$ws.Range("D46").Value = "'0.07713"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.88%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.04%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5532"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "

# Row 49
$ws.Range("E49").Value = "  +1.32%  "

# Row 50
$ws.Range("D50").Value = "'113.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "

# Row 51
$ws.Range("D51").Value = "'0.2931"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.92%  "

